# Swap row-pair data (everything except column A "id", and the constant
# columns C/D/E which are identical within each pair) between the two rows
# of each pair, per the commit's "base update" diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B and F:AC values are exchanged
$pairs = @(
    @(29, 30),
    @(42, 43),
    @(75, 76),
    @(90, 91),
    @(94, 95),
    @(96, 97),
    @(98, 99),
    @(102, 103),
    @(112, 113)
)

# Columns to swap: B (2) and F..AC (6..29)
$cols = @(2) + @(6..29)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($v1 -ne $v2) {
            $cell1.Value2 = $v2
            $cell2.Value2 = $v1
        }
    }
}
